# Add a set of "soft skills" rating columns (C:K) next to the existing
# username/password columns on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column titles, C1:K1 ---------------------
$headers = @(
    @{ Cell = "C1"; Text = "Dicipline" },
    @{ Cell = "D1"; Text = "Punctual, Time Management" },
    @{ Cell = "E1"; Text = "Work Ethic (Professionalism, Integrity, Honesty)" },
    @{ Cell = "F1"; Text = "Basic Communication" },
    @{ Cell = "G1"; Text = "Reliability" },
    @{ Cell = "H1"; Text = "Commitment" },
    @{ Cell = "I1"; Text = "Team work" },
    @{ Cell = "J1"; Text = "Listening" },
    @{ Cell = "K1"; Text = "Attention to details" }
)

foreach ($h in $headers) {
    $ws.Range($h.Cell).Value = $h.Text
    $ws.Range($h.Cell).Font.Size = 10
}

# F1 ("Basic Communication") is additionally highlighted with an
# Accent2 fill (Orange, Accent 2, Lighter 40%).
$ws.Range("F1").Interior.ThemeColor = 6
$ws.Range("F1").Interior.TintAndShade = 0.59999389629810485

# --- Data row (row 2): numeric ratings, C2:K2 --------------------------
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2

# Leave the active cell / selection on K1, matching the authored file.
$ws.Range("K1").Select()
